$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 857; everything below shifts down by one
# (old row 857 becomes row 858, ..., old row 929 becomes row 930).
$ws.Rows.Item(857).Insert()

$ws.Range("A857").Value = 3
$ws.Range("B857").Value = "Femacal de La Calera"
$ws.Range("C857").Value = "Coquimbo"
$ws.Range("D857").Value = 45106
$ws.Range("E857").Value = 5
$ws.Range("F857").Value = 100112045
$ws.Range("G857").Value = "Zapallo"
$ws.Range("H857").Value = "Camote"
$ws.Range("I857").Value = "1a (guarda)"
$ws.Range("J857").Value = 220
$ws.Range("K857").Value = 480
$ws.Range("L857").Value = 500
$ws.Range("M857").Value = 489
$ws.Range("N857").Value = "$/kilo (volumen en unidades)"
$ws.Range("O857").Value = "Provincia de Talca"
$ws.Range("P857").Value = 489
$ws.Range("Q857").Value = 1
$ws.Range("R857").Value = "Hortaliza"
